$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F3 2810 -> 2815, F4 253 -> 255
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2815
$ws1.Range("F4").Value = 255

# Sheet "全部类型" (sheet4): F4 2810 -> 2815, F6 253 -> 255
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2815
$ws4.Range("F6").Value = 255
